# Append daily weather observations for day-of-year rows 366-388
# (leap-year tail previously missing from the sheet) and refresh the
# active selection / scroll position to the top of the new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 366
$ws.Cells.Item(366, 6).Value = 13
$ws.Cells.Item(366, 7).Value = 6
$ws.Cells.Item(366, 8).Value = -2
$ws.Cells.Item(366, 9).Value = 5
$ws.Cells.Item(366, 10).Value = -1
$ws.Cells.Item(366, 11).Value = -4
$ws.Cells.Item(366, 12).Value = 87
$ws.Cells.Item(366, 13).Value = 63
$ws.Cells.Item(366, 14).Value = 33
$ws.Cells.Item(366, 15).Value = 1020
$ws.Cells.Item(366, 16).Value = 1016
$ws.Cells.Item(366, 17).Value = 1012
$ws.Cells.Item(366, 18).Value = 14
$ws.Cells.Item(366, 19).Value = 9
$ws.Cells.Item(366, 20).Value = 6
$ws.Cells.Item(366, 21).Value = 35
$ws.Cells.Item(366, 22).Value = 11
$ws.Cells.Item(366, 24).Value = 7.87
$ws.Cells.Item(366, 25).Value = 6
$ws.Cells.Item(366, 26).Value = "Rain"
$ws.Cells.Item(366, 27).Value = 182

# Row 367
$ws.Cells.Item(367, 6).Value = 18
$ws.Cells.Item(367, 7).Value = 12
$ws.Cells.Item(367, 8).Value = 6
$ws.Cells.Item(367, 9).Value = 5
$ws.Cells.Item(367, 10).Value = 2
$ws.Cells.Item(367, 11).Value = -3
$ws.Cells.Item(367, 12).Value = 93
$ws.Cells.Item(367, 13).Value = 57
$ws.Cells.Item(367, 14).Value = 18
$ws.Cells.Item(367, 15).Value = 1014
$ws.Cells.Item(367, 16).Value = 1012
$ws.Cells.Item(367, 17).Value = 1007
$ws.Cells.Item(367, 18).Value = 10
$ws.Cells.Item(367, 19).Value = 9
$ws.Cells.Item(367, 20).Value = 6
$ws.Cells.Item(367, 21).Value = 35
$ws.Cells.Item(367, 22).Value = 13
$ws.Cells.Item(367, 24).Value = 1.02
$ws.Cells.Item(367, 25).Value = 3
$ws.Cells.Item(367, 26).Value = "Rain"
$ws.Cells.Item(367, 27).Value = 184

# Row 368
$ws.Cells.Item(368, 6).Value = 13
$ws.Cells.Item(368, 7).Value = 8
$ws.Cells.Item(368, 8).Value = 3
$ws.Cells.Item(368, 9).Value = 3
$ws.Cells.Item(368, 10).Value = -1
$ws.Cells.Item(368, 11).Value = -7
$ws.Cells.Item(368, 12).Value = 93
$ws.Cells.Item(368, 13).Value = 60
$ws.Cells.Item(368, 14).Value = 15
$ws.Cells.Item(368, 15).Value = 1015
$ws.Cells.Item(368, 16).Value = 1012
$ws.Cells.Item(368, 17).Value = 1009
$ws.Cells.Item(368, 18).Value = 10
$ws.Cells.Item(368, 19).Value = 10
$ws.Cells.Item(368, 20).Value = 8
$ws.Cells.Item(368, 21).Value = 23
$ws.Cells.Item(368, 22).Value = 6
$ws.Cells.Item(368, 24).Value = 0
$ws.Cells.Item(368, 25).Value = 4
$ws.Cells.Item(368, 26).Value = "Rain"
$ws.Cells.Item(368, 27).Value = 202

# Row 369
$ws.Cells.Item(369, 6).Value = 14
$ws.Cells.Item(369, 7).Value = 8
$ws.Cells.Item(369, 8).Value = 1
$ws.Cells.Item(369, 9).Value = 2
$ws.Cells.Item(369, 10).Value = -1
$ws.Cells.Item(369, 11).Value = -7
$ws.Cells.Item(369, 12).Value = 93
$ws.Cells.Item(369, 13).Value = 63
$ws.Cells.Item(369, 14).Value = 16
$ws.Cells.Item(369, 15).Value = 1021
$ws.Cells.Item(369, 16).Value = 1017
$ws.Cells.Item(369, 17).Value = 1015
$ws.Cells.Item(369, 18).Value = 19
$ws.Cells.Item(369, 19).Value = 11
$ws.Cells.Item(369, 20).Value = 8
$ws.Cells.Item(369, 21).Value = 23
$ws.Cells.Item(369, 22).Value = 8
$ws.Cells.Item(369, 24).Value = 0.76
$ws.Cells.Item(369, 25).Value = 2
$ws.Cells.Item(369, 26).Value = "Rain"
$ws.Cells.Item(369, 27).Value = 209

# Row 370
$ws.Cells.Item(370, 6).Value = 18
$ws.Cells.Item(370, 7).Value = 9
$ws.Cells.Item(370, 8).Value = 0
$ws.Cells.Item(370, 9).Value = 2
$ws.Cells.Item(370, 10).Value = 0
$ws.Cells.Item(370, 11).Value = -2
$ws.Cells.Item(370, 12).Value = 93
$ws.Cells.Item(370, 13).Value = 57
$ws.Cells.Item(370, 14).Value = 19
$ws.Cells.Item(370, 15).Value = 1025
$ws.Cells.Item(370, 16).Value = 1022
$ws.Cells.Item(370, 17).Value = 1020
$ws.Cells.Item(370, 18).Value = 11
$ws.Cells.Item(370, 19).Value = 10
$ws.Cells.Item(370, 20).Value = 10
$ws.Cells.Item(370, 21).Value = 11
$ws.Cells.Item(370, 22).Value = 3
$ws.Cells.Item(370, 24).Value = 0
$ws.Cells.Item(370, 25).Value = 6
$ws.Cells.Item(370, 27).Value = 193

# Row 371
$ws.Cells.Item(371, 6).Value = 20
$ws.Cells.Item(371, 7).Value = 12
$ws.Cells.Item(371, 8).Value = 5
$ws.Cells.Item(371, 9).Value = 4
$ws.Cells.Item(371, 10).Value = 1
$ws.Cells.Item(371, 11).Value = -3
$ws.Cells.Item(371, 12).Value = 81
$ws.Cells.Item(371, 13).Value = 49
$ws.Cells.Item(371, 14).Value = 14
$ws.Cells.Item(371, 15).Value = 1023
$ws.Cells.Item(371, 16).Value = 1020
$ws.Cells.Item(371, 17).Value = 1014
$ws.Cells.Item(371, 18).Value = 10
$ws.Cells.Item(371, 19).Value = 10
$ws.Cells.Item(371, 20).Value = 10
$ws.Cells.Item(371, 21).Value = 32
$ws.Cells.Item(371, 22).Value = 10
$ws.Cells.Item(371, 24).Value = 0
$ws.Cells.Item(371, 25).Value = 6
$ws.Cells.Item(371, 27).Value = 188

# Row 372
$ws.Cells.Item(372, 6).Value = 13
$ws.Cells.Item(372, 7).Value = 11
$ws.Cells.Item(372, 8).Value = 9
$ws.Cells.Item(372, 9).Value = 6
$ws.Cells.Item(372, 10).Value = 3
$ws.Cells.Item(372, 11).Value = -2
$ws.Cells.Item(372, 12).Value = 82
$ws.Cells.Item(372, 13).Value = 55
$ws.Cells.Item(372, 14).Value = 24
$ws.Cells.Item(372, 15).Value = 1021
$ws.Cells.Item(372, 16).Value = 1018
$ws.Cells.Item(372, 17).Value = 1011
$ws.Cells.Item(372, 18).Value = 11
$ws.Cells.Item(372, 19).Value = 9
$ws.Cells.Item(372, 20).Value = 7
$ws.Cells.Item(372, 21).Value = 35
$ws.Cells.Item(372, 22).Value = 21
$ws.Cells.Item(372, 24).Value = 0
$ws.Cells.Item(372, 25).Value = 7
$ws.Cells.Item(372, 26).Value = "Rain"
$ws.Cells.Item(372, 27).Value = 188

# Row 373
$ws.Cells.Item(373, 6).Value = 11
$ws.Cells.Item(373, 7).Value = 9
$ws.Cells.Item(373, 8).Value = 8
$ws.Cells.Item(373, 9).Value = 8
$ws.Cells.Item(373, 10).Value = 6
$ws.Cells.Item(373, 11).Value = 4
$ws.Cells.Item(373, 12).Value = 93
$ws.Cells.Item(373, 13).Value = 75
$ws.Cells.Item(373, 14).Value = 55
$ws.Cells.Item(373, 15).Value = 1016
$ws.Cells.Item(373, 16).Value = 1011
$ws.Cells.Item(373, 17).Value = 1007
$ws.Cells.Item(373, 18).Value = 10
$ws.Cells.Item(373, 19).Value = 7
$ws.Cells.Item(373, 20).Value = 5
$ws.Cells.Item(373, 21).Value = 34
$ws.Cells.Item(373, 22).Value = 21
$ws.Cells.Item(373, 24).Value = 7.87
$ws.Cells.Item(373, 25).Value = 8
$ws.Cells.Item(373, 26).Value = "Rain"
$ws.Cells.Item(373, 27).Value = 186

# Row 374
$ws.Cells.Item(374, 6).Value = 11
$ws.Cells.Item(374, 7).Value = 7
$ws.Cells.Item(374, 8).Value = 4
$ws.Cells.Item(374, 9).Value = 7
$ws.Cells.Item(374, 10).Value = 5
$ws.Cells.Item(374, 11).Value = 3
$ws.Cells.Item(374, 12).Value = 93
$ws.Cells.Item(374, 13).Value = 82
$ws.Cells.Item(374, 14).Value = 61
$ws.Cells.Item(374, 15).Value = 1013
$ws.Cells.Item(374, 16).Value = 1010
$ws.Cells.Item(374, 17).Value = 1006
$ws.Cells.Item(374, 18).Value = 10
$ws.Cells.Item(374, 19).Value = 7
$ws.Cells.Item(374, 20).Value = 3
$ws.Cells.Item(374, 21).Value = 19
$ws.Cells.Item(374, 22).Value = 10
$ws.Cells.Item(374, 24).Value = 17.02
$ws.Cells.Item(374, 25).Value = 7
$ws.Cells.Item(374, 26).Value = "Rain"
$ws.Cells.Item(374, 27).Value = 182

# Row 375
$ws.Cells.Item(375, 6).Value = 13
$ws.Cells.Item(375, 7).Value = 8
$ws.Cells.Item(375, 8).Value = 3
$ws.Cells.Item(375, 9).Value = 4
$ws.Cells.Item(375, 10).Value = -1
$ws.Cells.Item(375, 11).Value = -5
$ws.Cells.Item(375, 12).Value = 100
$ws.Cells.Item(375, 13).Value = 57
$ws.Cells.Item(375, 14).Value = 18
$ws.Cells.Item(375, 15).Value = 1015
$ws.Cells.Item(375, 16).Value = 1011
$ws.Cells.Item(375, 17).Value = 1008
$ws.Cells.Item(375, 18).Value = 14
$ws.Cells.Item(375, 19).Value = 10
$ws.Cells.Item(375, 20).Value = 5
$ws.Cells.Item(375, 21).Value = 32
$ws.Cells.Item(375, 22).Value = 14
$ws.Cells.Item(375, 24).Value = 5.08
$ws.Cells.Item(375, 25).Value = 3
$ws.Cells.Item(375, 26).Value = "Rain"
$ws.Cells.Item(375, 27).Value = 184

# Row 376
$ws.Cells.Item(376, 6).Value = 13
$ws.Cells.Item(376, 7).Value = 7
$ws.Cells.Item(376, 8).Value = 2
$ws.Cells.Item(376, 9).Value = 1
$ws.Cells.Item(376, 10).Value = -2
$ws.Cells.Item(376, 11).Value = -6
$ws.Cells.Item(376, 12).Value = 93
$ws.Cells.Item(376, 13).Value = 55
$ws.Cells.Item(376, 14).Value = 17
$ws.Cells.Item(376, 15).Value = 1018
$ws.Cells.Item(376, 16).Value = 1017
$ws.Cells.Item(376, 17).Value = 1013
$ws.Cells.Item(376, 18).Value = 19
$ws.Cells.Item(376, 19).Value = 11
$ws.Cells.Item(376, 20).Value = 7
$ws.Cells.Item(376, 21).Value = 19
$ws.Cells.Item(376, 22).Value = 6
$ws.Cells.Item(376, 24).Value = 0.51
$ws.Cells.Item(376, 25).Value = 2
$ws.Cells.Item(376, 26).Value = "Rain"
$ws.Cells.Item(376, 27).Value = 175

# Row 377
$ws.Cells.Item(377, 6).Value = 14
$ws.Cells.Item(377, 7).Value = 6
$ws.Cells.Item(377, 8).Value = -2
$ws.Cells.Item(377, 9).Value = 0
$ws.Cells.Item(377, 10).Value = -3
$ws.Cells.Item(377, 11).Value = -8
$ws.Cells.Item(377, 12).Value = 93
$ws.Cells.Item(377, 13).Value = 49
$ws.Cells.Item(377, 14).Value = 12
$ws.Cells.Item(377, 15).Value = 1019
$ws.Cells.Item(377, 16).Value = 1016
$ws.Cells.Item(377, 17).Value = 1013
$ws.Cells.Item(377, 18).Value = 14
$ws.Cells.Item(377, 19).Value = 11
$ws.Cells.Item(377, 20).Value = 10
$ws.Cells.Item(377, 21).Value = 14
$ws.Cells.Item(377, 22).Value = 5
$ws.Cells.Item(377, 24).Value = 0
$ws.Cells.Item(377, 25).Value = 2
$ws.Cells.Item(377, 27).Value = 240

# Row 378
$ws.Cells.Item(378, 6).Value = 9
$ws.Cells.Item(378, 7).Value = 5
$ws.Cells.Item(378, 8).Value = 1
$ws.Cells.Item(378, 9).Value = 1
$ws.Cells.Item(378, 10).Value = -1
$ws.Cells.Item(378, 11).Value = -3
$ws.Cells.Item(378, 12).Value = 87
$ws.Cells.Item(378, 13).Value = 61
$ws.Cells.Item(378, 14).Value = 30
$ws.Cells.Item(378, 15).Value = 1022
$ws.Cells.Item(378, 16).Value = 1018
$ws.Cells.Item(378, 17).Value = 1014
$ws.Cells.Item(378, 18).Value = 14
$ws.Cells.Item(378, 19).Value = 10
$ws.Cells.Item(378, 20).Value = 9
$ws.Cells.Item(378, 21).Value = 23
$ws.Cells.Item(378, 22).Value = 10
$ws.Cells.Item(378, 24).Value = 0
$ws.Cells.Item(378, 25).Value = 5
$ws.Cells.Item(378, 26).Value = "Rain"
$ws.Cells.Item(378, 27).Value = 11

# Row 379
$ws.Cells.Item(379, 6).Value = 16
$ws.Cells.Item(379, 7).Value = 7
$ws.Cells.Item(379, 8).Value = -2
$ws.Cells.Item(379, 9).Value = -2
$ws.Cells.Item(379, 10).Value = -4
$ws.Cells.Item(379, 11).Value = -8
$ws.Cells.Item(379, 12).Value = 93
$ws.Cells.Item(379, 13).Value = 52
$ws.Cells.Item(379, 14).Value = 12
$ws.Cells.Item(379, 15).Value = 1024
$ws.Cells.Item(379, 16).Value = 1021
$ws.Cells.Item(379, 17).Value = 1017
$ws.Cells.Item(379, 18).Value = 19
$ws.Cells.Item(379, 19).Value = 12
$ws.Cells.Item(379, 20).Value = 10
$ws.Cells.Item(379, 21).Value = 11
$ws.Cells.Item(379, 22).Value = 3
$ws.Cells.Item(379, 24).Value = 0
$ws.Cells.Item(379, 25).Value = 1
$ws.Cells.Item(379, 27).Value = 148

# Row 380
$ws.Cells.Item(380, 6).Value = 17
$ws.Cells.Item(380, 7).Value = 8
$ws.Cells.Item(380, 8).Value = -1
$ws.Cells.Item(380, 9).Value = -1
$ws.Cells.Item(380, 10).Value = -4
$ws.Cells.Item(380, 11).Value = -9
$ws.Cells.Item(380, 12).Value = 93
$ws.Cells.Item(380, 13).Value = 47
$ws.Cells.Item(380, 14).Value = 8
$ws.Cells.Item(380, 15).Value = 1020
$ws.Cells.Item(380, 16).Value = 1018
$ws.Cells.Item(380, 17).Value = 1014
$ws.Cells.Item(380, 18).Value = 19
$ws.Cells.Item(380, 19).Value = 14
$ws.Cells.Item(380, 20).Value = 10
$ws.Cells.Item(380, 21).Value = 29
$ws.Cells.Item(380, 22).Value = 6
$ws.Cells.Item(380, 24).Value = 0
$ws.Cells.Item(380, 25).Value = 4
$ws.Cells.Item(380, 27).Value = 190

# Row 381
$ws.Cells.Item(381, 6).Value = 18
$ws.Cells.Item(381, 7).Value = 9
$ws.Cells.Item(381, 8).Value = 1
$ws.Cells.Item(381, 9).Value = 0
$ws.Cells.Item(381, 10).Value = -2
$ws.Cells.Item(381, 11).Value = -6
$ws.Cells.Item(381, 12).Value = 81
$ws.Cells.Item(381, 13).Value = 46
$ws.Cells.Item(381, 14).Value = 12
$ws.Cells.Item(381, 15).Value = 1018
$ws.Cells.Item(381, 16).Value = 1017
$ws.Cells.Item(381, 17).Value = 1013
$ws.Cells.Item(381, 18).Value = 19
$ws.Cells.Item(381, 19).Value = 12
$ws.Cells.Item(381, 20).Value = 10
$ws.Cells.Item(381, 21).Value = 29
$ws.Cells.Item(381, 22).Value = 10
$ws.Cells.Item(381, 24).Value = 0
$ws.Cells.Item(381, 25).Value = 2
$ws.Cells.Item(381, 27).Value = 189

# Row 382
$ws.Cells.Item(382, 6).Value = 20
$ws.Cells.Item(382, 7).Value = 10
$ws.Cells.Item(382, 8).Value = 1
$ws.Cells.Item(382, 9).Value = 1
$ws.Cells.Item(382, 10).Value = -2
$ws.Cells.Item(382, 11).Value = -6
$ws.Cells.Item(382, 12).Value = 81
$ws.Cells.Item(382, 13).Value = 46
$ws.Cells.Item(382, 14).Value = 11
$ws.Cells.Item(382, 15).Value = 1017
$ws.Cells.Item(382, 16).Value = 1015
$ws.Cells.Item(382, 17).Value = 1012
$ws.Cells.Item(382, 18).Value = 14
$ws.Cells.Item(382, 19).Value = 13
$ws.Cells.Item(382, 20).Value = 10
$ws.Cells.Item(382, 21).Value = 23
$ws.Cells.Item(382, 22).Value = 5
$ws.Cells.Item(382, 24).Value = 0
$ws.Cells.Item(382, 25).Value = 1
$ws.Cells.Item(382, 27).Value = 241

# Row 383
$ws.Cells.Item(383, 6).Value = 22
$ws.Cells.Item(383, 7).Value = 12
$ws.Cells.Item(383, 8).Value = 2
$ws.Cells.Item(383, 9).Value = 1
$ws.Cells.Item(383, 10).Value = -1
$ws.Cells.Item(383, 11).Value = -9
$ws.Cells.Item(383, 12).Value = 81
$ws.Cells.Item(383, 13).Value = 44
$ws.Cells.Item(383, 14).Value = 7
$ws.Cells.Item(383, 15).Value = 1018
$ws.Cells.Item(383, 16).Value = 1016
$ws.Cells.Item(383, 17).Value = 1011
$ws.Cells.Item(383, 18).Value = 14
$ws.Cells.Item(383, 19).Value = 12
$ws.Cells.Item(383, 20).Value = 10
$ws.Cells.Item(383, 21).Value = 19
$ws.Cells.Item(383, 22).Value = 3
$ws.Cells.Item(383, 24).Value = 0
$ws.Cells.Item(383, 25).Value = 2
$ws.Cells.Item(383, 27).Value = 220

# Row 384
$ws.Cells.Item(384, 6).Value = 22
$ws.Cells.Item(384, 7).Value = 12
$ws.Cells.Item(384, 8).Value = 1
$ws.Cells.Item(384, 9).Value = 3
$ws.Cells.Item(384, 10).Value = 0
$ws.Cells.Item(384, 11).Value = -8
$ws.Cells.Item(384, 12).Value = 87
$ws.Cells.Item(384, 13).Value = 50
$ws.Cells.Item(384, 14).Value = 10
$ws.Cells.Item(384, 15).Value = 1019
$ws.Cells.Item(384, 16).Value = 1017
$ws.Cells.Item(384, 17).Value = 1012
$ws.Cells.Item(384, 18).Value = 19
$ws.Cells.Item(384, 19).Value = 12
$ws.Cells.Item(384, 20).Value = 10
$ws.Cells.Item(384, 21).Value = 23
$ws.Cells.Item(384, 22).Value = 3
$ws.Cells.Item(384, 24).Value = 0
$ws.Cells.Item(384, 25).Value = 2
$ws.Cells.Item(384, 27).Value = 184

# Row 385
$ws.Cells.Item(385, 6).Value = 22
$ws.Cells.Item(385, 7).Value = 12
$ws.Cells.Item(385, 8).Value = 3
$ws.Cells.Item(385, 9).Value = 2
$ws.Cells.Item(385, 10).Value = 1
$ws.Cells.Item(385, 11).Value = -3
$ws.Cells.Item(385, 12).Value = 81
$ws.Cells.Item(385, 13).Value = 47
$ws.Cells.Item(385, 14).Value = 12
$ws.Cells.Item(385, 15).Value = 1020
$ws.Cells.Item(385, 16).Value = 1017
$ws.Cells.Item(385, 17).Value = 1012
$ws.Cells.Item(385, 18).Value = 14
$ws.Cells.Item(385, 19).Value = 11
$ws.Cells.Item(385, 20).Value = 10
$ws.Cells.Item(385, 21).Value = 26
$ws.Cells.Item(385, 22).Value = 5
$ws.Cells.Item(385, 24).Value = 0
$ws.Cells.Item(385, 25).Value = 2
$ws.Cells.Item(385, 27).Value = 64

# Row 386
$ws.Cells.Item(386, 6).Value = 22
$ws.Cells.Item(386, 7).Value = 13
$ws.Cells.Item(386, 8).Value = 4
$ws.Cells.Item(386, 9).Value = 4
$ws.Cells.Item(386, 10).Value = 0
$ws.Cells.Item(386, 11).Value = -29
$ws.Cells.Item(386, 12).Value = 87
$ws.Cells.Item(386, 13).Value = 49
$ws.Cells.Item(386, 14).Value = 4
$ws.Cells.Item(386, 15).Value = 1019
$ws.Cells.Item(386, 16).Value = 1017
$ws.Cells.Item(386, 17).Value = 1011
$ws.Cells.Item(386, 18).Value = 19
$ws.Cells.Item(386, 19).Value = 11
$ws.Cells.Item(386, 20).Value = 10
$ws.Cells.Item(386, 21).Value = 11
$ws.Cells.Item(386, 22).Value = 3
$ws.Cells.Item(386, 24).Value = 0
$ws.Cells.Item(386, 25).Value = 2
$ws.Cells.Item(386, 27).Value = 154

# Row 387
$ws.Cells.Item(387, 6).Value = 21
$ws.Cells.Item(387, 7).Value = 13
$ws.Cells.Item(387, 8).Value = 6
$ws.Cells.Item(387, 9).Value = 8
$ws.Cells.Item(387, 10).Value = 4
$ws.Cells.Item(387, 11).Value = 2
$ws.Cells.Item(387, 12).Value = 77
$ws.Cells.Item(387, 13).Value = 51
$ws.Cells.Item(387, 14).Value = 22
$ws.Cells.Item(387, 15).Value = 1020
$ws.Cells.Item(387, 16).Value = 1017
$ws.Cells.Item(387, 17).Value = 1012
$ws.Cells.Item(387, 18).Value = 14
$ws.Cells.Item(387, 19).Value = 10
$ws.Cells.Item(387, 20).Value = 8
$ws.Cells.Item(387, 21).Value = 26
$ws.Cells.Item(387, 22).Value = 5
$ws.Cells.Item(387, 24).Value = 0
$ws.Cells.Item(387, 25).Value = 5
$ws.Cells.Item(387, 26).Value = "Rain"
$ws.Cells.Item(387, 27).Value = 128

# Row 388
$ws.Cells.Item(388, 6).Value = 15
$ws.Cells.Item(388, 7).Value = 11
$ws.Cells.Item(388, 8).Value = 8
$ws.Cells.Item(388, 9).Value = 11
$ws.Cells.Item(388, 10).Value = 8
$ws.Cells.Item(388, 11).Value = 7
$ws.Cells.Item(388, 12).Value = 94
$ws.Cells.Item(388, 13).Value = 83
$ws.Cells.Item(388, 14).Value = 53
$ws.Cells.Item(388, 15).Value = 1019
$ws.Cells.Item(388, 16).Value = 1016
$ws.Cells.Item(388, 17).Value = 1011
$ws.Cells.Item(388, 18).Value = 10
$ws.Cells.Item(388, 19).Value = 8
$ws.Cells.Item(388, 20).Value = 6
$ws.Cells.Item(388, 21).Value = 29
$ws.Cells.Item(388, 22).Value = 5
$ws.Cells.Item(388, 24).Value = 5.08
$ws.Cells.Item(388, 25).Value = 6
$ws.Cells.Item(388, 26).Value = "Rain"
$ws.Cells.Item(388, 27).Value = 210

# Excel recomputes the dimension/sheetView automatically; just move
# the selection to match where the author left the cursor afterward.
$ws.Range("F1:AA388").Select()
